# Workbook edit: add a "Player Info" sheet ahead of "ODI Batting", and
# replace the match-card URL column on "ODI Batting" with a bare match code.

$wb = $excel.ActiveWorkbook
$odiSheet = $wb.Worksheets.Item("ODI Batting")

# --- 1. Update "ODI Batting": MATCH_CARD_LINK -> MATCH_CODE -----------------
$odiSheet.Range("D1").Value = "MATCH_CODE"
$odiSheet.Range("D2:D3").NumberFormat = "@"
$odiSheet.Range("D2").Value = "4694"
$odiSheet.Range("D3").Value = "4696"

# --- 2. Insert a new "Player Info" sheet before "ODI Batting" ---------------
$newSheet = $wb.Worksheets.Add($odiSheet)
$newSheet.Name = "Player Info"

# Header row, styled to match the workbook's existing header look
# (bold font, thin box border, centered horizontally, top-aligned vertically)
$headerRange = $newSheet.Range("A1:D1")
$headerRange.Borders.LineStyle = 1
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

$newSheet.Range("A1").Value = "ID"
$newSheet.Range("B1").Value = "NAME"
$newSheet.Range("C1").Value = "BATTING_HAND"
$newSheet.Range("D1").Value = "BOWL_STYLE"

# Data row
$newSheet.Range("A2").NumberFormat = "@"
$newSheet.Range("A2").Value = "7122"
$newSheet.Range("B2").Value = "Murray Commins"
$newSheet.Range("C2").Value = "Left Handed"
$newSheet.Range("D2").Value = "Does Not Bowl | Unknown"
